$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range (A1:G43) so removed rows/cells do not linger
$ws.Range("A1:I43").ClearContents()

# --- Column widths: widen column A, keep B-G, add H and I ---
$ws.Columns("A:A").ColumnWidth = 46.833333333333336
$ws.Columns("H:H").ColumnWidth = 36.833333333333336
$ws.Columns("I:I").ColumnWidth = 45.833333333333336

# --- Table data ---
$ws.Range("A1").Value = "Brand name"
$ws.Range("B1").Value = "Generic name"
$ws.Range("C1").Value = "Doses"

$ws.Range("A2").Value = "Bijuve/Bijuva"
$ws.Range("B2").Value = "Estradiol"
$ws.Range("C2").Value = "Estradiol 1mg / Progesterone 100mg capsules"

$ws.Range("A3").Value = "Blissel"
$ws.Range("B3").Value = "Estriol"
$ws.Range("C3").Value = "Estriol 50micrograms/g vaginal gel with applicator"

$ws.Range("A4").Value = "Elleste Duet"
$ws.Range("C4").Value = "Generic Elleste Duet 1mg tablets"
$ws.Range("D4").Value = "Generic Elleste Duet 2mg tablets"

$ws.Range("A5").Value = "Elleste Duet Conti"
$ws.Range("B5").Value = "Estradiol"
$ws.Range("C5").Value = "Estradiol 2mg / Norethisterone acetate 1mg tablets"

$ws.Range("A6").Value = "Elleste Solo"
$ws.Range("B6").Value = "Estradiol"
$ws.Range("C6").Value = "Estradiol 1mg tablets"
$ws.Range("D6").Value = "Estradiol 2mg tablets"

$ws.Range("A7").Value = "Estraderm MX"
$ws.Range("B7").Value = "Estradiol"
$ws.Range("C7").Value = "Estradiol 25micrograms/24hours transdermal patches"
$ws.Range("D7").Value = "Estradiol 50micrograms/24hours transdermal patches"
$ws.Range("E7").Value = "Estradiol 75micrograms/24hours transdermal patches"
$ws.Range("F7").Value = "Estradiol 100micrograms/24hours transdermal patches"

$ws.Range("A8").Value = "Estradot"
$ws.Range("B8").Value = "Estradiol"
$ws.Range("C8").Value = "Estradiol 25micrograms/24hours transdermal patches"
$ws.Range("D8").Value = "Estradiol 37.5micrograms/24hours transdermal patches"
$ws.Range("E8").Value = "Estradiol 50micrograms/24hours transdermal patches"
$ws.Range("F8").Value = "Estradiol 75micrograms/24hours transdermal patches"
$ws.Range("G8").Value = "Estradiol 100micrograms/24hours transdermal patches"

$ws.Range("A9").Value = "Estring (90-day preparation)"
$ws.Range("B9").Value = "Estradiol"
$ws.Range("C9").Value = "Estradiol 7.5micrograms/24hours vaginal delivery system"

$ws.Range("A10").Value = "Evorel"
$ws.Range("B10").Value = "Estradiol"
$ws.Range("C10").Value = "Estradiol 25micrograms/24hours transdermal patches"
$ws.Range("D10").Value = "Estradiol 50micrograms/24hours transdermal patches"
$ws.Range("E10").Value = "Estradiol 75micrograms/24hours transdermal patches"
$ws.Range("F10").Value = "Estradiol 100micrograms/24hours transdermal patches"

$ws.Range("A11").Value = "Evorel Conti"
$ws.Range("B11").Value = "Estradiol"
$ws.Range("C11").Value = "Estradiol 50micrograms/24hours / Norethisterone 170micrograms/24hours transdermal patches"

$ws.Range("A12").Value = "Evorel Sequi"
$ws.Range("C12").Value = "Generic Evorel Sequi transdermal patches"

$ws.Range("A13").Value = "Femoston"
$ws.Range("C13").Value = "Generic Femoston 1/10mg tablets"
$ws.Range("D13").Value = "Generic Femoston 2/10mg tablets"

$ws.Range("A14").Value = "Femoston Conti"
$ws.Range("B14").Value = "Estradiol"
$ws.Range("C14").Value = "Estradiol 500micrograms / Dydrogesterone 2.5mg tablets"
$ws.Range("D14").Value = "Estradiol 1mg / Dydrogesterone 5mg tablets"

$ws.Range("A15").Value = "Femseven"
$ws.Range("B15").Value = "Estradiol"
$ws.Range("C15").Value = "Estradiol 50micrograms/24hours transdermal patches"
$ws.Range("D15").Value = "Estradiol 75micrograms/24hours transdermal patches"
$ws.Range("E15").Value = "Estradiol 100micrograms/24hours transdermal patches"

$ws.Range("A16").Value = "Femseven Conti"
$ws.Range("B16").Value = "Estradiol"
$ws.Range("C16").Value = "Estradiol 50micrograms/24hours / Levonorgestrel 7micrograms/24hours transdermal patches"

$ws.Range("A17").Value = "FemSeven Sequi Phase 1"
$ws.Range("B17").Value = "Estradiol"
$ws.Range("C17").Value = "Estradiol 50micrograms/24hours transdermal patches"

$ws.Range("A18").Value = "FemSeven Sequi Phase 2"
$ws.Range("B18").Value = "Estradiol"
$ws.Range("C18").Value = "Estradiol 50micrograms/24hours / Levonorgestrel 10micrograms/24hours transdermal patches"

$ws.Range("A19").Value = "Gepretix"
$ws.Range("B19").Value = "Progesterone"
$ws.Range("C19").Value = "Progesterone micronised 100mg capsules"
$ws.Range("D19").Value = "Progesterone micronised 200mg capsules"

$ws.Range("A20").Value = "Gina"
$ws.Range("B20").Value = "Estradiol"
$ws.Range("C20").Value = "Estradiol 10microgram pessaries"

$ws.Range("A21").Value = "Imvaggis"
$ws.Range("B21").Value = "Estriol"
$ws.Range("C21").Value = "Estriol 30microgram pessaries"

$ws.Range("A22").Value = "Indivina"
$ws.Range("C22").Value = "Estradiol valerate 1mg / Medroxyprogesterone 2.5mg tablets"
$ws.Range("D22").Value = "Estradiol valerate 1mg / Medroxyprogesterone 5mg tablets"
$ws.Range("E22").Value = "Estradiol valerate 2mg / Medroxyprogesterone 5mg table"

$ws.Range("A23").Value = "Intrarosa"
$ws.Range("B23").Value = "Prasterone"
$ws.Range("C23").Value = "Prasterone 6.5mg pessaries"

$ws.Range("A24").Value = "Kliofem"
$ws.Range("B24").Value = "Estriadiol"
$ws.Range("C24").Value = "Estradiol 2mg / Norethisterone acetate 1mg tablets"

$ws.Range("A25").Value = "Kliovance"
$ws.Range("B25").Value = "Estriadiol"
$ws.Range("C25").Value = "Estradiol 1mg / Norethisterone acetate 500microgram tablets"

$ws.Range("A26").Value = "Lenzetto"
$ws.Range("B26").Value = "Estriadiol"
$ws.Range("C26").Value = "Estradiol 1.53mg/dose transdermal spray"

$ws.Range("A27").Value = "Livial"
$ws.Range("B27").Value = "Tibolone"
$ws.Range("C27").Value = "Tibolone 2.5mg tablets"

$ws.Range("A28").Value = "Mirena"
$ws.Range("C28").Value = "Levonorgestrel 20micrograms/24hours intrauterine device"

$ws.Range("A29").Value = "Nalvee"
$ws.Range("C29").Value = "Dydrogesterone 10mg tablets"

$ws.Range("A30").Value = "Non-branded products not listed elsewhere in the table"
$ws.Range("C30").Value = "Conjugated oestrogens 300microgram tablets"
$ws.Range("D30").Value = "Conjugated oestrogens 625microgram tablets"
$ws.Range("E30").Value = "Conjugated oestrogens 1.25microgram tablets"
$ws.Range("F30").Value = "Estriol 500microgram pessaries"
$ws.Range("G30").Value = "Estriol 0.01% vaginal cream with applicator"
$ws.Range("H30").Value = "Estriol 0.01% vaginal cream"
$ws.Range("I30").Value = "Medroxyprogesterone 1.5mg modified-release tablets"

$ws.Range("A31").Value = "Novofem"
$ws.Range("C31").Value = "Generic Novofem tablets"

$ws.Range("A32").Value = "Oestrogel/dose"
$ws.Range("B32").Value = "Estradiol"
$ws.Range("C32").Value = "Estradiol 0.06% gel (750microgram per actuation)"

$ws.Range("A33").Value = "Progynova"
$ws.Range("B33").Value = "Estriadiol"
$ws.Range("C33").Value = "Estradiol valerate 1mg tablets"
$ws.Range("D33").Value = "Estradiol valerate 2mg tablets"

$ws.Range("A34").Value = "Progynova TS"
$ws.Range("B34").Value = "Estriadiol"
$ws.Range("C34").Value = "Estradiol 50micrograms/24hours transdermal patches"
$ws.Range("D34").Value = "Estradiol 100micrograms/24hours transdermal patches"

$ws.Range("A35").Value = "Sandrena"
$ws.Range("B35").Value = "Estriadiol"
$ws.Range("C35").Value = "Estradiol 500microgram gel sachets"
$ws.Range("D35").Value = "Estradiol 1mg gel sachets"

$ws.Range("A36").Value = "Tridestra"
$ws.Range("C36").Value = "Generic Tridestra tablets"

$ws.Range("A37").Value = "Trisequens"
$ws.Range("C37").Value = "Generic Trisequens tablets"

$ws.Range("A38").Value = "Utrogestan"
$ws.Range("B38").Value = "Progesterone"
$ws.Range("C38").Value = "Progesterone micronised 100mg capsules"

$ws.Range("A39").Value = "Vagifem"
$ws.Range("B39").Value = "Estradiol"
$ws.Range("C39").Value = "Estradiol 10microgram pessaries"

$ws.Range("A40").Value = "Vagirux"
$ws.Range("B40").Value = "Estradiol"
$ws.Range("C40").Value = "Estradiol 10microgram pessaries"

$ws.Range("A41").Value = "Zumenon"
$ws.Range("B41").Value = "Estradiol"
$ws.Range("C41").Value = "Estradiol 1mg tablets"
$ws.Range("D41").Value = "Estradiol 2mg tablets"

# --- Restore selection to match the saved view state ---
$ws.Range("D29").Select()